$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new flight entry for Friday, Jan 13 (Wizz Air / A321)
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(18, 3).Value = "10:10 AM"
$ws.Cells.Item(18, 4).Value = "W92065"
$ws.Cells.Item(18, 5).Value = "London"
$ws.Cells.Item(18, 6).Value = "(LTN)"
$ws.Cells.Item(18, 7).Value = "Wizz Air "
$ws.Cells.Item(18, 8).Value = "A321"
$ws.Cells.Item(18, 9).Value = "(G-WUKJ)"
$ws.Cells.Item(18, 10).Value = "10:02 AM"
$ws.Cells.Item(18, 11).Borders.LineStyle = -4142
$ws.Cells.Item(18, 12).Value = "0 hours, -8 minutes"
$ws.Cells.Item(18, 13).Borders.LineStyle = -4142

# Row 19: new flight entry for Friday, Jan 13 (Ryanair / B738)
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(19, 3).Value = "1:30 PM"
$ws.Cells.Item(19, 4).Value = "FR6640"
$ws.Cells.Item(19, 5).Value = "London"
$ws.Cells.Item(19, 6).Value = "(LTN)"
$ws.Cells.Item(19, 7).Value = "Ryanair "
$ws.Cells.Item(19, 8).Value = "B738"
$ws.Cells.Item(19, 9).Value = "(EI-DPL)"
$ws.Cells.Item(19, 10).Value = "1:24 PM"
$ws.Cells.Item(19, 11).Borders.LineStyle = -4142
$ws.Cells.Item(19, 12).Value = "0 hours, -6 minutes"
$ws.Cells.Item(19, 13).Borders.LineStyle = -4142
